$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Normal style: turn off "overflow punctuation" (w:overflowPunct)
#    true -> false. In this object model this boolean is surfaced as
#    ParagraphFormat.HangingPunctuation.
# ------------------------------------------------------------------
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.HangingPunctuation = $false

# ------------------------------------------------------------------
# 2) Add the new "ListLabel NN" character styles (43 through 61) that
#    round out the list-numbering character styles already present
#    (ListLabel1 .. ListLabel42). Each one only carries a qFormat flag
#    plus a small rPr (complex-script font, occasionally a size or
#    bold flag).
# ------------------------------------------------------------------

# id, complex-script font name (or $null), point size (or $null), bold flag
$specs = @(
    @(43, "Symbol",      11,   $false),
    @(44, "Courier New", $null,$false),
    @(45, "Wingdings",   $null,$false),
    @(46, "Symbol",      $null,$false),
    @(47, "Courier New", $null,$false),
    @(48, "Wingdings",   $null,$false),
    @(49, "Symbol",      $null,$false),
    @(50, "Courier New", $null,$false),
    @(51, "Wingdings",   $null,$false),
    @(52, $null,         $null,$true),
    @(53, "Symbol",      $null,$false),
    @(54, "Courier New", $null,$false),
    @(55, "Wingdings",   $null,$false),
    @(56, "Symbol",      $null,$false),
    @(57, "Courier New", $null,$false),
    @(58, "Wingdings",   $null,$false),
    @(59, "Symbol",      $null,$false),
    @(60, "Courier New", $null,$false),
    @(61, "Wingdings",   $null,$false)
)

foreach ($spec in $specs) {
    $num  = $spec[0]
    $font = $spec[1]
    $size = $spec[2]
    $bold = $spec[3]

    $styleId = "ListLabel$num"
    $s = $d.Styles.Add($styleId, 2)
    $s.NameLocal = "ListLabel $num"
    $s.QuickStyle = $true

    if ($font) {
        $s.Font.NameBi = $font
    }
    if ($size) {
        $s.Font.Size = $size
    }
    if ($bold) {
        $s.Font.Bold = $true
    }
}

Write-Output "done"
